# Project Assignment #2 - "conceptual classes" bullet list was reworked to
# match the new class diagram: several items were renamed/reordered, some
# were dropped, and some new ones were introduced. We rebuild the whole
# bulleted list (same ListParagraph/numId list it already used) in one shot
# via Range.InsertXML, which lets us control the exact <w:proofErr>/
# <w:bookmarkStart|End> markup the new content needs, the same way real
# Word's spell/grammar checker and "_GoBack" bookmark tracking would.

$d = $word.ActiveDocument

# --- locate the bullet list to replace ---------------------------------
# It starts at the "SMS Messenger" item and runs to the very end of the
# document body (the blank trailing ListParagraph after "Screen" is part of
# the block being replaced too).
$count = $d.Paragraphs.Count
$idxFirst = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd()
    if ($t -eq "SMS Messenger") { $idxFirst = $i }
}

$firstPara = $d.Paragraphs($idxFirst)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$rng = $d.Range($firstPara.Range.Start, $lastPara.Range.End)

# --- build the replacement list items -----------------------------------
# Each entry: display text, proofErr marks wrapping the run (if any), and
# whether the "_GoBack" bookmark (the last edit position) sits on it.
$items = @(
    @{ Text = "android.app.Activity";                Proof = "spellgram" },
    @{ Text = "android.fragment.app.DialogFragment";  Proof = "spellgram" },
    @{ Text = "ErrorDialog";                          Proof = "spell" },
    @{ Text = "WarningDialog";                        Proof = "spell" },
    @{ Text = "MessageActivity";                      Proof = "spell" },
    @{ Text = "ConversationActivity";                 Proof = "spell" },
    @{ Text = "Conversation";                         Proof = "none" },
    @{ Text = "Message";                              Proof = "none" },
    @{ Text = "ReceiverService";                       Proof = "spell" },
    @{ Text = "android.app.Service";                  Proof = "spellgram" },
    @{ Text = "Blacklist";                             Proof = "none" },
    @{ Text = "Contact";                               Proof = "none" },
    @{ Text = "Contact Manager";                       Proof = "none"; Bookmark = $true }
)

$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

$bodyXml = ""
foreach ($item in $items) {
    $run = "<w:r><w:t>" + $item.Text + "</w:t></w:r>"

    switch ($item.Proof) {
        "spellgram" {
            $run = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + $run + '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>'
        }
        "spell" {
            $run = '<w:proofErr w:type="spellStart"/>' + $run + '<w:proofErr w:type="spellEnd"/>'
        }
        default {
        }
    }

    if ($item.Bookmark -eq $true) {
        $run = $run + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    }

    $bodyXml = $bodyXml + "<w:p>" + $pPr + $run + "</w:p>"
}

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- apply ---------------------------------------------------------------
[void]$rng.InsertXML($packageXml)
